$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

# Fix header text for column N (remove trailing space: "Correction " -> "Correction")
$ws.Range("N1").Value = "Correction"

# Add new header "Serviced by " in column O, matching the bold/centered header style of N1
$ws.Range("O1").Value = "Serviced by "
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

# Column N (rows 2-12) becomes "nan" like its row siblings (was blank before)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}

# Create the new column O cells (rows 2-12) as blank cells so the used range
# and sheet dimension extend to column O, without introducing a new style
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Borders.LineStyle = 0
}
